# Applies the symbol-list refresh described in the commit message:
# coin prices/volumes updated, and the GateToken..HotbitToken block
# rotated by one row (HotbitToken moved from row 19 to row 24).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store numeric-looking values as
# plain text in this workbook (e.g. "330.08", "1.25%"). Mark the cells
# we are about to touch as Text first so Excel does not silently
# convert them to numbers/percentages and drop significant digits.
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "E25",
    "D26",
    "E26",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "E41",
    "D42",
    "E42",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "E46",
    "D47",
    "D48",
    "E48",
    "D49",
    "E49",
    "E50",
    "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2: BNB
$ws.Range("D2").Value = "330.08"
$ws.Range("E2").Value = "1.25%"

# Row 3: OKB
$ws.Range("D3").Value = "44.19"
$ws.Range("E3").Value = "-0.88%"

# Row 4: HuobiToken
$ws.Range("D4").Value = "5.493"
$ws.Range("E4").Value = "-1.85%"

# Row 5: Cronos
$ws.Range("D5").Value = "0.08005"
$ws.Range("E5").Value = "-0.77%"

# Row 6: FTXToken
$ws.Range("D6").Value = "1.975"
$ws.Range("E6").Value = "3.74%"

# Row 7: BTSEToken
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "2.579"
$ws.Range("E7").Value = "-4.61%"

# Row 8: MXToken
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9513"
$ws.Range("E8").Value = "0.61%"

# Row 9: LiechtensteinCryptoassetsExchange
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1099"
$ws.Range("E9").Value = "-5.63%"

# Row 10: WazirX
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1908"
$ws.Range("E10").Value = "2.28%"

# Row 11: MCDex
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").Value = "10.51"
$ws.Range("E11").Value = "24.17%"

# Row 12: MandalaExchangeToken
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09947"
$ws.Range("E12").Value = "-0.11%"

# Row 13: BitrueCoin
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.04789"
$ws.Range("E13").Value = "12.87%"

# Row 14: BitMartToken
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.1064"
$ws.Range("E14").Value = "-0.19%"

# Row 15: BitForexToken
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001274"
$ws.Range("E15").Value = "-0.90%"

# Row 16: CoinExToken
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "0.04086"
$ws.Range("E16").Value = "-2.92%"

# Row 17: TigerCash
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.005945"
$ws.Range("E17").Value = "0.78%"

# Row 18: LEO
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.370"
$ws.Range("E18").Value = "-6.24%"

# Row 19: GateToken
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "4.391"
$ws.Range("E19").Value = "1.92%"

# Row 20: BitpandaEcosystemToken
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3491"
$ws.Range("E20").Value = "-0.22%"

# Row 21: ProBitToken
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1420"
$ws.Range("E21").Value = "3.54%"

# Row 22: ZBToken
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "0.2587"
$ws.Range("E22").Value = "-0.89%"

# Row 23: BitKan
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "0.001273"
$ws.Range("E23").Value = "2.51%"

# Row 24: HotbitToken
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "0.004374"
$ws.Range("E24").Value = "-2.06%"

# Row 25: NitroEx
$ws.Range("E25").Value = "1.64%"

# Row 26: UpBots
$ws.Range("D26").Value = "0.0003745"
$ws.Range("E26").Value = "-6.21%"

# Row 38: One
$ws.Range("D38").Value = "0.02586"
$ws.Range("E38").Value = "-1.64%"

# Row 39: IDEX
$ws.Range("D39").Value = "0.05694"
$ws.Range("E39").Value = "4.67%"

# Row 40: KickToken
$ws.Range("D40").Value = "0.007563"
$ws.Range("E40").Value = "-1.69%"

# Row 41: BKEXToken
$ws.Range("E41").Value = "0.23%"

# Row 42: Dexo
$ws.Range("D42").Value = "0.007359"
$ws.Range("E42").Value = "2.98%"

# Row 43: CEJI
$ws.Range("E43").Value = "-0.49%"

# Row 44: LocalTraders
$ws.Range("D44").Value = "0.008354"
$ws.Range("E44").Value = "-2.58%"

# Row 45: CoinLion
$ws.Range("D45").Value = "0.00007132"
$ws.Range("E45").Value = "-0.09%"

# Row 46: Kangarootoken
$ws.Range("E46").Value = "-0.06%"

# Row 47: ACDXExchange
$ws.Range("D47").Value = "0.0005803"

# Row 48: CoinbaseStockToken
$ws.Range("D48").Value = "0.003531"
$ws.Range("E48").Value = "55.40%"

# Row 49: BOLO
$ws.Range("D49").Value = "0.003556"
$ws.Range("E49").Value = "0.08%"

# Row 50: CryptobidCoin
$ws.Range("E50").Value = "-0.06%"

# Row 51: SpecialPowerGold
$ws.Range("E51").Value = "-0.06%"
